$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (formerly "Strike#"). Regenerated values per row.
$gValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    25 = 3
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
